# urbs_intertemporal_2050 / 2045.xlsx
# "Implemented timeseries to supim file, demand file next"
#
# This adds 11 more timestep rows (t=2..12) to the SupIm sheet, mirroring
# the existing t=1 row's values, makes SupIm the active sheet/selection,
# and consolidates the two conditional-formatting rules on the Process
# sheet that used to straddle row 11 into a single contiguous rule.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) SupIm sheet: extend the timeseries from t=1 (row 3) down to t=12
#    (row 14), copying values/format from the existing t=1 row.
# ---------------------------------------------------------------------
$supIm = $wb.Worksheets.Item("SupIm")

for ($row = 4; $row -le 14; $row++) {
    $t = $row - 2
    $supIm.Cells.Item($row, 1).Value = $t
    $supIm.Cells.Item($row, 2).Value = 0.481
    $supIm.Cells.Item($row, 3).Value = 0.3
    $supIm.Cells.Item($row, 4).Value = 0.207
}

# Copy the formatting of the template row (row 3) down onto the new rows
# (values were already written above so they keep being numeric).
$supIm.Range("A3:D3").Copy() | Out-Null
$supIm.Range("A4:D14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Sheet becomes the active tab, with a fresh selection (matches the
# author scrolling down and clicking on J18 after entering the data).
$supIm.Activate() | Out-Null
$supIm.Range("J18").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Process sheet: merge the two conditional formatting rules that
#    used to cover (A12:C13, A11, C11) and (B11) separately into a
#    single rule covering the contiguous block A11:C13.
# ---------------------------------------------------------------------
$process = $wb.Worksheets.Item("Process")

$mainRule = $process.Cells.FormatConditions.Item(1)
$mainRule.ModifyAppliesToRange($process.Range("A11:C13"))
$process.Cells.FormatConditions.Item(2).Delete()
$process.Cells.FormatConditions.Item(1).Priority = 1
